# daily auto push: 2026-01-17 06:45 UTC
# Insert a new observation row for 2026/01/17 (time 13, ranking 201) at
# row 667, pushing the existing rows 667-708 down to 668-709.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 667; everything below shifts down by one.
$ws.Rows(667).Insert()

# Date column: enter as literal text (leading apostrophe stops Excel from
# auto-converting the "yyyy/mm/dd"-looking text into a real date), then
# reset the style back to Normal so no extra number-format/quote-prefix
# flag is left behind on the cell.
$ws.Cells.Item(667, 1).Value = "'2026/01/17"
$ws.Cells.Item(667, 1).Style = "Normal"

# Day-of-week, hour, and ranking columns.
$ws.Cells.Item(667, 2).Value = "土"
$ws.Cells.Item(667, 3).Value = 13
$ws.Cells.Item(667, 4).Value = 201
